# Apply the "want-to-go" (F column) count bumps across all four sheets, and
# insert the new duplicate event row on "全部类型" (pushing rows 37-47 down
# to 38-48), matching the upstream data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "展览" (Exhibitions) sheet - simple numeric bumps in column F
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 18
$wsExpo.Range("F3").Value  = 2739
$wsExpo.Range("F5").Value  = 19775
$wsExpo.Range("F6").Value  = 78
$wsExpo.Range("F7").Value  = 2281
$wsExpo.Range("F8").Value  = 752
$wsExpo.Range("F13").Value = 254
$wsExpo.Range("F15").Value = 374
$wsExpo.Range("F16").Value = 75

# ---------------------------------------------------------------------
# 2) "演出" (Performances) sheet - simple numeric bumps in column F
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value  = 198
$wsShow.Range("F7").Value  = 290
$wsShow.Range("F8").Value  = 132
$wsShow.Range("F15").Value = 82

# ---------------------------------------------------------------------
# 3) "本地生活" (Local life) sheet - simple numeric bumps in column F
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 648
$wsLocal.Range("F4").Value = 594

# ---------------------------------------------------------------------
# 4) "全部类型" (All types) sheet - numeric bumps in column F ...
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 648
$wsAll.Range("F4").Value  = 594
$wsAll.Range("F5").Value  = 198
$wsAll.Range("F6").Value  = 18
$wsAll.Range("F8").Value  = 2739
$wsAll.Range("F10").Value = 19775
$wsAll.Range("F13").Value = 78
$wsAll.Range("F15").Value = 290
$wsAll.Range("F16").Value = 2281
$wsAll.Range("F17").Value = 752
$wsAll.Range("F18").Value = 132
$wsAll.Range("F23").Value = 254
$wsAll.Range("F28").Value = 374
$wsAll.Range("F29").Value = 75
$wsAll.Range("F36").Value = 82

# ... plus a brand-new row inserted right after row 36: a duplicate of the
# (now F=82) "平田雄也&小池亮介" event, pushing the old rows 37-47 down to
# 38-48 and growing the sheet dimension from A1:I47 to A1:I48.
$wsAll.Rows.Item(37).Insert()

# The freshly inserted row inherits a slightly-off style from Insert(); pull
# the real formatting back in from the row beneath it (which still has the
# original, untouched formatting) before writing any values.
$wsAll.Range("A38:I38").Copy()
$wsAll.Range("A37:I37").PasteSpecial(-4122)  # xlPasteFormats

# Column B holds dates as plain text (e.g. "2024-11-10"); force Text format
# before assignment so it is not silently reinterpreted as a date serial,
# then paste the sibling row's format back in so no stray numFmt sticks.
$wsAll.Range("B37").NumberFormat = "@"
$wsAll.Range("B37").Value = "2024-11-10"
$wsAll.Range("B38").Copy()
$wsAll.Range("B37").PasteSpecial(-4122)  # xlPasteFormats

$wsAll.Range("A37").Value = 36
$wsAll.Range("C37").Value = "广州·平田雄也&小池亮介2024粉丝见面会"
$wsAll.Range("D37").Value = "金花街道中山七路333号1906科技圆区3号楼109-1铺、110-1铺、111-1铺 音乐唐人馆"
$wsAll.Range("E37").Value = "2024.11.10 13:00-11.10 18:00"
$wsAll.Range("F37").Value = 82
$wsAll.Range("G37").Value = 480
$wsAll.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=92655"
$wsAll.Range("I37").Value = "//i2.hdslb.com/bfs/openplatform/202409/UkhOeOwe1726658317935.jpeg"

Write-Output "edits applied"
